# Trade #31 closed at 2026-02-17 08:28:16 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet updates ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.37
$wsSummary.Range("B4").Value = -0.63
$wsSummary.Range("B6").Value = 31
$wsSummary.Range("B8").Value = 18
$wsSummary.Range("B9").Value = 22.58

# --- Strategy Status sheet updates (MarketMaking row, row 4) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.37
$wsStatus.Range("D4").Value = 31
$wsStatus.Range("E4").Value = -0.63
$wsStatus.Range("F4").Value = -0.63
$wsStatus.Range("G4").Value = 22.58

# --- Append new trade row (row 32) to "All Trades" and "MarketMaking" sheets ---
$newRow = @{
    A = 31
    B = "2026-02-17"
    C = "08:28:10"
    D = "MarketMaking"
    E = "UP"
    F = 0.98
    G = 0.97
    H = "CLOSED"
    I = -1.0204
    J = -0.01
    K = 99.37
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.14
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A32").Value = $newRow.A
    # B32 looks like a date (yyyy-mm-dd); force text entry so Excel
    # doesn't auto-convert it to a date serial number, then restore the
    # cell's plain "Normal" style (no lingering number-format override).
    $ws.Range("B32").NumberFormat = "@"
    $ws.Range("B32").Value = $newRow.B
    $ws.Range("B32").Style = "Normal"
    $ws.Range("C32").Value = $newRow.C
    $ws.Range("D32").Value = $newRow.D
    $ws.Range("E32").Value = $newRow.E
    $ws.Range("F32").Value = $newRow.F
    $ws.Range("G32").Value = $newRow.G
    $ws.Range("H32").Value = $newRow.H
    $ws.Range("I32").Value = $newRow.I
    $ws.Range("J32").Value = $newRow.J
    $ws.Range("K32").Value = $newRow.K
    $ws.Range("L32").Value = $newRow.L
    $ws.Range("M32").Value = $newRow.M
    $ws.Range("N32").Value = $newRow.N
    $ws.Range("O32").Value = $newRow.O
    $ws.Range("P32").Value = $newRow.P
    $ws.Range("Q32").Value = $newRow.Q
}
